# edit.ps1 -- applies the text-content edits described by the commit diff.
#
# The diff contains a large number of purely cosmetic / re-save artifacts
# (XML namespace prefix changes on <w:document>, dropped wp14:/w14:/w15:
# markup-compatibility extensions, <w:tblLook> attribute simplification,
# <w:lang> tweak, <a14:useLocalDpi> extra xmlns declarations, regenerated
# <w:proofErr> spellcheck/grammar-check run-splitting, latentStyles table
# rewrite, theme panose removal, etc.) that are artifacts of the document
# having been re-saved by a different Word build/version and are not
# reachable (or meaningful) through the Word object model -- they do not
# change the document's visible text or content.
#
# The only edits that change the actual document content are both in the
# last body paragraphs of word/document.xml:
#
#   1. A new sentence is appended at the end of the paragraph that ends
#      "...y un falso (0) para filtrarlo.":
#         " //programadores mayores de 30 años"
#
#   2. In the "Nota 0" paragraph, the space between "de las" and
#      "funciones" is removed, producing "de lasfunciones".
#
# Both are applied below using Find/Replace on the document's Content
# range, matching Word's own COM automation semantics.

$d = $word.ActiveDocument

# 1) Append the new comment sentence right after "...para filtrarlo."
$rng = $d.Content
$found1 = $rng.Find.Execute("filtrarlo.")
if ($found1) {
    $rng.Collapse(0)
    $rng.InsertAfter(" //programadores mayores de 30 años")
}

# 2) Close the gap between "de las" and "funciones" in the Nota 0 paragraph.
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("de las funciones", $false, $false, $false, $false, $false, $true, 1, $false, "de lasfunciones", 2)
